# Update countries & provincias Spain
# - Refresh the "last updated" timestamp
# - Update case/death/recovered statistics for a handful of countries
# - Lebanon ("Libano") moved up past Malaui/Nicaragua in the ranking after
#   its numbers were refreshed: remove its stale row and re-insert it (with
#   the new figures) above Malaui, pushing Malaui/Nicaragua down one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp header ---
$ws.Range("A1").Value = "Datos actualizados a 25 de Julio de 2020 a las 20:33"

# --- Numeric updates for existing rows (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Estados Unidos (row 4)
$ws.Range("B4").Value = 4283058
$ws.Range("C4").Value = 34731
$ws.Range("D4").Value = 2039742
$ws.Range("E4").Value = 2094325
$ws.Range("G4").Value = 501
$ws.Range("H4").Value = 148991

# Brasil (row 5)
$ws.Range("B5").Value = 2355920
$ws.Range("C5").Value = 7720
$ws.Range("E5").Value = 678077
$ws.Range("G5").Value = 177
$ws.Range("H5").Value = 85562

# India (row 6)
$ws.Range("B6").Value = 1385494
$ws.Range("C6").Value = 48472
$ws.Range("D6").Value = 886235
$ws.Range("E6").Value = 467163
$ws.Range("G6").Value = 690
$ws.Range("H6").Value = 32096

# Alemania (row 21)
$ws.Range("B21").Value = 206203
$ws.Range("C21").Value = 243
$ws.Range("E21").Value = 6602

# Canada (row 24)
$ws.Range("B24").Value = 113515
$ws.Range("C24").Value = 309
$ws.Range("D24").Value = 99111

# Israel (row 41)
$ws.Range("B41").Value = 60496
$ws.Range("C41").Value = 1021
$ws.Range("D41").Value = 26882
$ws.Range("E41").Value = 33159
$ws.Range("G41").Value = 7
$ws.Range("H41").Value = 455

# Marruecos (row 66)
$ws.Range("B66").Value = 19645
$ws.Range("C66").Value = 811
$ws.Range("D66").Value = 16282
$ws.Range("E66").Value = 3058
$ws.Range("G66").Value = 6
$ws.Range("H66").Value = 305

# Sri Lanka (row 114)
$ws.Range("B114").Value = 2770
$ws.Range("C114").Value = 6
$ws.Range("E114").Value = 656

# --- Reorder Libano ---
# Row 107 currently holds Libano's stale data; delete it so Tailandia's
# block (row 108 onward) shifts up by one.
$ws.Rows.Item(107).Delete()

# Insert a fresh row above Malaui (currently row 105) and push Malaui /
# Nicaragua back down to rows 106 / 107.
$ws.Rows.Item(105).Insert()

$ws.Range("A105").Value = "Libano"
$ws.Range("B105").Value = 3582
$ws.Range("C105").Value = 175
$ws.Range("D105").Value = 1671
$ws.Range("E105").Value = 1864
$ws.Range("F105").Value = 0
$ws.Range("G105").Value = 1
$ws.Range("H105").Value = 47
